# Rename "Sheet2" to "offlinechat" and populate it with a small
# name/email/message header row, then make it the active sheet
# with A2 selected (mirrors the "offline chat" feature commit).

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "offlinechat"

$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "email"
$ws2.Range("C1").Value = "message"

$ws2.Activate()
$ws2.Range("A2").Select()
